$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds price strings; several new values look numeric
# (e.g. "1.001") and Excel would silently coerce them to floats,
# losing the original text formatting - force Text format first.
$textCells = @('D5', 'D6', 'D7', 'D8', 'D9', 'D11', 'D12', 'D14', 'D15', 'D17', 'D18', 'D20', 'D21', 'D23', 'D24', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.142.39'
$ws.Range('E2').Value = '  -0.48%  '

$ws.Range('D3').Value = '1.885.21'
$ws.Range('E3').Value = '  -0.82%  '

$ws.Range('E4').Value = '  +0.38%  '

$ws.Range('D5').Value = '322.13'
$ws.Range('E5').Value = '  -2.88%  '

$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.39%  '

$ws.Range('D7').Value = '0.4703'
$ws.Range('E7').Value = '  +2.01%  '

$ws.Range('D8').Value = '0.4025'
$ws.Range('E8').Value = '  -2.74%  '

$ws.Range('D9').Value = '47.23'
$ws.Range('E9').Value = '  -0.94%  '

$ws.Range('E10').Value = '  -0.14%  '

$ws.Range('D11').Value = '0.9920'
$ws.Range('E11').Value = '  -1.67%  '

$ws.Range('D12').Value = '22.40'
$ws.Range('E12').Value = '  +0.47%  '

$ws.Range('D13').Value = '1.878.70'
$ws.Range('E13').Value = '  -2.52%  '

$ws.Range('D14').Value = '5.859'
$ws.Range('E14').Value = '  -1.52%  '

$ws.Range('D15').Value = '7.014'
$ws.Range('E15').Value = '  -1.37%  '

$ws.Range('E16').Value = '  +0.37%  '

$ws.Range('D17').Value = '88.71'
$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').Value = '0.06616'
$ws.Range('E18').Value = '  +0.77%  '

$ws.Range('E19').Value = '  -0.43%  '

$ws.Range('D20').Value = '17.38'
$ws.Range('E20').Value = '  -1.46%  '

$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.60%  '

$ws.Range('D22').Value = '29.149.38'
$ws.Range('E22').Value = '  -0.22%  '

$ws.Range('D23').Value = '5.480'
$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').Value = '11.60'
$ws.Range('E24').Value = '  +1.41%  '

$ws.Range('E25').Value = '  -0.83%  '

$ws.Range('D26').Value = '2.119.42'
$ws.Range('E26').Value = '  -1.38%  '

$ws.Range('D27').Value = '155.06'
$ws.Range('E27').Value = '  -0.90%  '

$ws.Range('D28').Value = '19.59'
$ws.Range('E28').Value = '  -0.62%  '

$ws.Range('D29').Value = '6.033'
$ws.Range('E29').Value = '  +7.34%  '

$ws.Range('D30').Value = '2.068'
$ws.Range('E30').Value = '  -2.17%  '

$ws.Range('D31').Value = '117.41'
$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('D32').Value = '1.034'
$ws.Range('E32').Value = '  -0.85%  '

$ws.Range('D33').Value = '0.09439'
$ws.Range('E33').Value = '  +0.59%  '

$ws.Range('D34').Value = '3.540'
$ws.Range('E34').Value = '  +0.44%  '

$ws.Range('D35').Value = '1.380'
$ws.Range('E35').Value = '  -3.18%  '

$ws.Range('D36').Value = '5.343'
$ws.Range('E36').Value = '  -0.19%  '

$ws.Range('D37').Value = '0.06061'
$ws.Range('E37').Value = '  -0.28%  '

$ws.Range('E38').Value = '  -0.71%  '

$ws.Range('D39').Value = '1.172'
$ws.Range('E39').Value = '  -0.33%  '

$ws.Range('D40').Value = '8.010'
$ws.Range('E40').Value = '  -5.09%  '

$ws.Range('D41').Value = '0.5797'
$ws.Range('E41').Value = '  -0.57%  '

$ws.Range('D42').Value = '0.1821'
$ws.Range('E42').Value = '  -0.47%  '

$ws.Range('D43').Value = '2.443'
$ws.Range('E43').Value = '  +5.13%  '

$ws.Range('D44').Value = '9.980'
$ws.Range('E44').Value = '  -1.93%  '

$ws.Range('D45').Value = '1.269'
$ws.Range('E45').Value = '  +1.84%  '

$ws.Range('D46').Value = '0.07664'
$ws.Range('E46').Value = '  +1.68%  '

$ws.Range('D47').Value = '12.03'
$ws.Range('E47').Value = '  -0.90%  '

$ws.Range('D48').Value = '0.5460'
$ws.Range('E48').Value = '  -1.03%  '

$ws.Range('D49').Value = '1.897'
$ws.Range('E49').Value = '  -1.40%  '

$ws.Range('D50').Value = '113.63'
$ws.Range('E50').Value = '  +1.52%  '

$ws.Range('D51').Value = '44.36'
$ws.Range('E51').Value = '  +0.95%  '
